$d = $word.ActiveDocument
$searchStart = 0

$rr0 = $d.Content
$rr0.Start = $searchStart
$rr0.End = $d.Content.End
$ok0 = $rr0.Find.Execute("550840.064", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok0) { Write-Host "NOT FOUND: 550840.064 (#0)" }
$xml0 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="08C475A5" w14:textId="085EAA8D" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>550840,064</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr0.InsertXML($xml0)
$searchStart = $rr0.End

$rr1 = $d.Content
$rr1.Start = $searchStart
$rr1.End = $d.Content.End
$ok1 = $rr1.Find.Execute("43470.752", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok1) { Write-Host "NOT FOUND: 43470.752 (#1)" }
$xml1 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="73E1E62C" w14:textId="2E92D51D" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>42517,477</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr1.InsertXML($xml1)
$searchStart = $rr1.End

$rr2 = $d.Content
$rr2.Start = $searchStart
$rr2.End = $d.Content.End
$ok2 = $rr2.Find.Execute("550840.064", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok2) { Write-Host "NOT FOUND: 550840.064 (#2)" }
$xml2 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="37410563" w14:textId="46AE72D6" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>550840,064</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr2.InsertXML($xml2)
$searchStart = $rr2.End

$rr3 = $d.Content
$rr3.Start = $searchStart
$rr3.End = $d.Content.End
$ok3 = $rr3.Find.Execute("42604.126", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok3) { Write-Host "NOT FOUND: 42604.126 (#3)" }
$xml3 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="54787676" w14:textId="0AB36FFE" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>43451,868</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr3.InsertXML($xml3)
$searchStart = $rr3.End

$rr4 = $d.Content
$rr4.Start = $searchStart
$rr4.End = $d.Content.End
$ok4 = $rr4.Find.Execute("550197.729", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok4) { Write-Host "NOT FOUND: 550197.729 (#4)" }
$xml4 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2BCF73A5" w14:textId="740503C8" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>550197,729</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr4.InsertXML($xml4)
$searchStart = $rr4.End

$rr5 = $d.Content
$rr5.Start = $searchStart
$rr5.End = $d.Content.End
$ok5 = $rr5.Find.Execute("43483.396", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok5) { Write-Host "NOT FOUND: 43483.396 (#5)" }
$xml5 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="50943FF5" w14:textId="09BE0428" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>43483,396</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr5.InsertXML($xml5)
$searchStart = $rr5.End

$rr6 = $d.Content
$rr6.Start = $searchStart
$rr6.End = $d.Content.End
$ok6 = $rr6.Find.Execute("550210.799", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok6) { Write-Host "NOT FOUND: 550210.799 (#6)" }
$xml6 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="230A087A" w14:textId="49067B1C" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>550210</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>799</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr6.InsertXML($xml6)
$searchStart = $rr6.End

$rr7 = $d.Content
$rr7.Start = $searchStart
$rr7.End = $d.Content.End
$ok7 = $rr7.Find.Execute("42479.668", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok7) { Write-Host "NOT FOUND: 42479.668 (#7)" }
$xml7 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="52F3F16F" w14:textId="2C68614E" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>42479</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>668</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr7.InsertXML($xml7)
$searchStart = $rr7.End

$rr8 = $d.Content
$rr8.Start = $searchStart
$rr8.End = $d.Content.End
$ok8 = $rr8.Find.Execute("550208.291", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok8) { Write-Host "NOT FOUND: 550208.291 (#8)" }
$xml8 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="284D1391" w14:textId="39275B0F" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>550208</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>291</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr8.InsertXML($xml8)
$searchStart = $rr8.End

$rr9 = $d.Content
$rr9.Start = $searchStart
$rr9.End = $d.Content.End
$ok9 = $rr9.Find.Execute("42329.905", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok9) { Write-Host "NOT FOUND: 42329.905 (#9)" }
$xml9 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2FF45946" w14:textId="224DAE1A" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>42329</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>905</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr9.InsertXML($xml9)
$searchStart = $rr9.End

$rr10 = $d.Content
$rr10.Start = $searchStart
$rr10.End = $d.Content.End
$ok10 = $rr10.Find.Execute("550855.932", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok10) { Write-Host "NOT FOUND: 550855.932 (#10)" }
$xml10 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="61039883" w14:textId="4319BD65" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>550855</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>932</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr10.InsertXML($xml10)
$searchStart = $rr10.End

$rr11 = $d.Content
$rr11.Start = $searchStart
$rr11.End = $d.Content.End
$ok11 = $rr11.Find.Execute("42847.353", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok11) { Write-Host "NOT FOUND: 42847.353 (#11)" }
$xml11 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="61FEC9D5" w14:textId="5CC1AA34" w:rsidR="0043769A" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="0043769A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>42847</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>353</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr11.InsertXML($xml11)
$searchStart = $rr11.End

$rr12 = $d.Content
$rr12.Start = $searchStart
$rr12.End = $d.Content.End
$ok12 = $rr12.Find.Execute("541049.338", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok12) { Write-Host "NOT FOUND: 541049.338 (#12)" }
$xml12 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="06D7E5CF" w14:textId="568D4A6B" w:rsidR="00BB3168" w:rsidRPr="000B34DE" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>541049</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>338</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr12.InsertXML($xml12)
$searchStart = $rr12.End

$rr13 = $d.Content
$rr13.Start = $searchStart
$rr13.End = $d.Content.End
$ok13 = $rr13.Find.Execute("64624.054", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok13) { Write-Host "NOT FOUND: 64624.054 (#13)" }
$xml13 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2397F4B1" w14:textId="056CBF2B" w:rsidR="00BB3168" w:rsidRPr="000B34DE" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>64624</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>054</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr13.InsertXML($xml13)
$searchStart = $rr13.End

$rr14 = $d.Content
$rr14.Start = $searchStart
$rr14.End = $d.Content.End
$ok14 = $rr14.Find.Execute("541049.338", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok14) { Write-Host "NOT FOUND: 541049.338 (#14)" }
$xml14 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="511D01C9" w14:textId="37EE955C" w:rsidR="00BB3168" w:rsidRPr="000B34DE" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>541049</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>338</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr14.InsertXML($xml14)
$searchStart = $rr14.End

$rr15 = $d.Content
$rr15.Start = $searchStart
$rr15.End = $d.Content.End
$ok15 = $rr15.Find.Execute("62826.344", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok15) { Write-Host "NOT FOUND: 62826.344 (#15)" }
$xml15 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="60DA5B40" w14:textId="56E4C1DD" w:rsidR="00BB3168" w:rsidRPr="000B34DE" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>62826</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>344</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr15.InsertXML($xml15)
$searchStart = $rr15.End

$rr16 = $d.Content
$rr16.Start = $searchStart
$rr16.End = $d.Content.End
$ok16 = $rr16.Find.Execute("541049.272", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok16) { Write-Host "NOT FOUND: 541049.272 (#16)" }
$xml16 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4F4F577F" w14:textId="20EE1427" w:rsidR="00BB3168" w:rsidRPr="000B34DE" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>541049</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>272</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr16.InsertXML($xml16)
$searchStart = $rr16.End

$rr17 = $d.Content
$rr17.Start = $searchStart
$rr17.End = $d.Content.End
$ok17 = $rr17.Find.Execute("63101.728", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok17) { Write-Host "NOT FOUND: 63101.728 (#17)" }
$xml17 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="31DC7DE0" w14:textId="6887DFFA" w:rsidR="00BB3168" w:rsidRPr="000B34DE" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>63101</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>728</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr17.InsertXML($xml17)
$searchStart = $rr17.End

$rr18 = $d.Content
$rr18.Start = $searchStart
$rr18.End = $d.Content.End
$ok18 = $rr18.Find.Execute("541064.705", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok18) { Write-Host "NOT FOUND: 541064.705 (#18)" }
$xml18 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="08F5D3A2" w14:textId="0883A66B" w:rsidR="00BB3168" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>541064</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>705</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr18.InsertXML($xml18)
$searchStart = $rr18.End

$rr19 = $d.Content
$rr19.Start = $searchStart
$rr19.End = $d.Content.End
$ok19 = $rr19.Find.Execute("67797.867", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok19) { Write-Host "NOT FOUND: 67797.867 (#19)" }
$xml19 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="59C49A6F" w14:textId="3D940C38" w:rsidR="00BB3168" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>67797</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>867</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr19.InsertXML($xml19)
$searchStart = $rr19.End

$rr20 = $d.Content
$rr20.Start = $searchStart
$rr20.End = $d.Content.End
$ok20 = $rr20.Find.Execute("541064.650", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok20) { Write-Host "NOT FOUND: 541064.650 (#20)" }
$xml20 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6DDD2F95" w14:textId="73A5A6C5" w:rsidR="00BB3168" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>541064</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>650</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr20.InsertXML($xml20)
$searchStart = $rr20.End

$rr21 = $d.Content
$rr21.Start = $searchStart
$rr21.End = $d.Content.End
$ok21 = $rr21.Find.Execute("63628.854", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok21) { Write-Host "NOT FOUND: 63628.854 (#21)" }
$xml21 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="59776A54" w14:textId="34FBB905" w:rsidR="00BB3168" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>63628</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>854</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr21.InsertXML($xml21)
$searchStart = $rr21.End

$rr22 = $d.Content
$rr22.Start = $searchStart
$rr22.End = $d.Content.End
$ok22 = $rr22.Find.Execute("541064.646", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok22) { Write-Host "NOT FOUND: 541064.646 (#22)" }
$xml22 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="777E9A42" w14:textId="37F5CC8F" w:rsidR="00BB3168" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>541064</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>646</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr22.InsertXML($xml22)
$searchStart = $rr22.End

$rr23 = $d.Content
$rr23.Start = $searchStart
$rr23.End = $d.Content.End
$ok23 = $rr23.Find.Execute("63037.975", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok23) { Write-Host "NOT FOUND: 63037.975 (#23)" }
$xml23 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4EF0D50E" w14:textId="6BB55852" w:rsidR="00BB3168" w:rsidRPr="0043769A" w:rsidRDefault="00BB3168" w:rsidP="00BB3168"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>63037</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof w:val="0"/><w:color w:val="000000"/></w:rPr><w:t>975</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rr23.InsertXML($xml23)
$searchStart = $rr23.End
